$d = $word.ActiveDocument

# Locate the paragraph that ends with "...avoidant/restrictive food intake
# disorder;" (it carries the two bookmarkEnd elements); the new paragraph
# is inserted right after it, before the "having experienced ... sexual
# abuse ..." paragraph.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*avoidant/restrictive food intake disorder;*") {
        $anchorIndex = $i
        break
    }
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)

$r = $anchorPara.Range
$r.Collapse(0)            # wdCollapseEnd - move to the end of the paragraph
$r.InsertParagraphAfter() # splits in a new, empty paragraph right after it

$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = "another para."

# Match the formatting of the anchor paragraph: style LV2, indent 1418,
# explicit "no numbering" override (ilvl 0 / numId 0).
$newPara.Style = $anchorPara.Style
$newPara.Range.ParagraphFormat.LeftIndent = $anchorPara.Range.ParagraphFormat.LeftIndent
$newPara.Range.ListFormat.RemoveNumbers()
